$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = New-Object 'object[,]' 132,1
$newValues[0,0] = 9.891639143228531
$newValues[1,0] = 10.36087130308151
$newValues[2,0] = 11.8361004948616
$newValues[3,0] = 14.00294080972671
$newValues[4,0] = 16.68157880306244
$newValues[5,0] = 19.6415261387825
$newValues[6,0] = 21.28996098041534
$newValues[7,0] = 20.12264679670334
$newValues[8,0] = 17.33410066366196
$newValues[9,0] = 13.82261809110642
$newValues[10,0] = 11.3167012989521
$newValues[11,0] = 10.2355012357235
$newValues[12,0] = 9.879003047943115
$newValues[13,0] = 10.67941368222237
$newValues[14,0] = 12.1199037194252
$newValues[15,0] = 13.99530737400055
$newValues[16,0] = 16.21080392599106
$newValues[17,0] = 17.93669948577881
$newValues[18,0] = 18.54963983297348
$newValues[19,0] = 17.42080775499344
$newValues[20,0] = 15.31175798177719
$newValues[21,0] = 13.10456585884094
$newValues[22,0] = 11.64428151845932
$newValues[23,0] = 10.76767483949661
$newValues[24,0] = 10.35594524741173
$newValues[25,0] = 11.01092126965523
$newValues[26,0] = 12.23676578998566
$newValues[27,0] = 13.8316814661026
$newValues[28,0] = 15.30383940935135
$newValues[29,0] = 16.54835612773895
$newValues[30,0] = 16.68902498483658
$newValues[31,0] = 15.93063229322433
$newValues[32,0] = 14.5379891872406
$newValues[33,0] = 13.10927699804306
$newValues[34,0] = 12.00802515745163
$newValues[35,0] = 11.31026089787483
$newValues[36,0] = 10.99008926749229
$newValues[37,0] = 11.37078860998154
$newValues[38,0] = 12.32919547557831
$newValues[39,0] = 13.47526555657387
$newValues[40,0] = 14.61593506336212
$newValues[41,0] = 15.45086355209351
$newValues[42,0] = 15.52021702528
$newValues[43,0] = 15.08116254210472
$newValues[44,0] = 14.20994863510132
$newValues[45,0] = 13.23060195446014
$newValues[46,0] = 12.38990760445595
$newValues[47,0] = 11.78912883400917
$newValues[48,0] = 11.52648814320564
$newValues[49,0] = 11.73472375869751
$newValues[50,0] = 12.35189292430878
$newValues[51,0] = 13.21577271819115
$newValues[52,0] = 14.07049551010132
$newValues[53,0] = 14.66752423644066
$newValues[54,0] = 14.77699332237244
$newValues[55,0] = 14.54809026122093
$newValues[56,0] = 14.01846586465836
$newValues[57,0] = 13.34612061977387
$newValues[58,0] = 12.70343384742737
$newValues[59,0] = 12.21338245272636
$newValues[60,0] = 11.96641364097595
$newValues[61,0] = 12.04066054224968
$newValues[62,0] = 12.4623919069767
$newValues[63,0] = 13.06974229812622
$newValues[64,0] = 13.67847652435303
$newValues[65,0] = 14.11930069923401
$newValues[66,0] = 14.26013411283493
$newValues[67,0] = 14.1597319483757
$newValues[68,0] = 13.85013732910156
$newValues[69,0] = 13.3988609790802
$newValues[70,0] = 12.92597382068634
$newValues[71,0] = 12.53632600903511
$newValues[72,0] = 12.3083705663681
$newValues[73,0] = 12.29182554483414
$newValues[74,0] = 12.56682641506195
$newValues[75,0] = 12.98061012625694
$newValues[76,0] = 13.41162049174309
$newValues[77,0] = 13.7455205321312
$newValues[78,0] = 13.89438174962998
$newValues[79,0] = 13.87261844277382
$newValues[80,0] = 13.70437644124031
$newValues[81,0] = 13.41258300542831
$newValues[82,0] = 13.07547482252121
$newValues[83,0] = 12.77330368757248
$newValues[84,0] = 12.56985723376274
$newValues[85,0] = 12.51929370164871
$newValues[86,0] = 12.66268625259399
$newValues[87,0] = 12.93439031839371
$newValues[88,0] = 13.23629688620567
$newValues[89,0] = 13.48900035023689
$newValues[90,0] = 13.62997987866402
$newValues[91,0] = 13.65535491704941
$newValues[92,0] = 13.57781972289085
$newValues[93,0] = 13.39948019385338
$newValues[94,0] = 13.16799955368042
$newValues[95,0] = 12.94052075743675
$newValues[96,0] = 12.77031826376915
$newValues[97,0] = 12.70114566087723
$newValues[98,0] = 12.74653133749962
$newValues[99,0] = 12.91738070845604
$newValues[100,0] = 13.12456869482994
$newValues[101,0] = 13.31404345035553
$newValues[102,0] = 13.43940713405609
$newValues[103,0] = 13.49265109300613
$newValues[104,0] = 13.47234467864036
$newValues[105,0] = 13.37295964360237
$newValues[106,0] = 13.22127188444138
$newValues[107,0] = 13.05642740726471
$newValues[108,0] = 12.91924899220467
$newValues[109,0] = 12.84126339554787
$newValues[110,0] = 12.84355016350746
$newValues[111,0] = 12.91958661675453
$newValues[112,0] = 13.05487121343613
$newValues[113,0] = 13.19326181411743
$newValues[114,0] = 13.3001391351223
$newValues[115,0] = 13.36885353326798
$newValues[116,0] = 13.38385654687881
$newValues[117,0] = 13.33870493769646
$newValues[118,0] = 13.24669657349586
$newValues[119,0] = 13.13288305401802
$newValues[120,0] = 13.02604970932007
$newValues[121,0] = 12.95219788551331
$newValues[122,0] = 12.9273505628109
$newValues[123,0] = 12.95628342628479
$newValues[124,0] = 13.02703846693039
$newValues[125,0] = 13.11689837574959
$newValues[126,0] = 13.20564469099045
$newValues[127,0] = 13.27475913167
$newValues[128,0] = 13.30789882540703
$newValues[129,0] = 13.29776228666305
$newValues[130,0] = 13.24890460968018
$newValues[131,0] = 13.17557340860367

$range = $ws.Range("F2:F133")
$range.Value = $newValues
